# Insert a new data row at row 306 (pushing the existing rows 306:373 down
# to 307:374) and populate the new row with the latest weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("306:306").Insert()

$ws.Cells.Item(306, 1).Value  = 10
$ws.Cells.Item(306, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(306, 3).Value  = "La Araucanía"
$ws.Cells.Item(306, 4).Value  = 44889
$ws.Cells.Item(306, 5).Value  = 9
$ws.Cells.Item(306, 6).Value  = 100114013
$ws.Cells.Item(306, 7).Value  = "Zanahoria"
$ws.Cells.Item(306, 8).Value  = "Sin especificar"
$ws.Cells.Item(306, 9).Value  = "Primera"
$ws.Cells.Item(306, 10).Value = 375
$ws.Cells.Item(306, 11).Value = 10000
$ws.Cells.Item(306, 12).Value = 12000
$ws.Cells.Item(306, 13).Value = 11333
$ws.Cells.Item(306, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(306, 15).Value = "Región del Bíobío"
$ws.Cells.Item(306, 16).Value = 567
$ws.Cells.Item(306, 17).Value = 20
$ws.Cells.Item(306, 18).Value = "Hortaliza"
